# Regenerate save_data: use K (strikeouts) column values instead of the
# previous Strike# values. Only column G (header "K") on rows 2-18 changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value, taken from the recalculated s_vals
$kValues = [ordered]@{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 6
    7  = 3
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 8
    14 = 1
    15 = 3
    16 = 4
    17 = 3
    18 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$wb.Save()
